# Append "  (This is a change – Version for main branch)" to the end of
# the first paragraph's text ("This is a Microsoft word document."),
# matching the target OOXML: the original run keeps two trailing spaces,
# and the parenthetical note is added as three separate red-colored runs.

$d = $word.ActiveDocument

# First paragraph currently reads: "This is a Microsoft word document."
$p1 = $d.Paragraphs(1).Range
$end = $p1.End - 1   # position right after the final period, before the paragraph mark

# 1) Append two trailing spaces to the existing (black) run.
$rSpace = $d.Range($end, $end)
$rSpace.InsertAfter("  ")
$pos = $end + 2

# 2) Append the red parenthetical text as three separate runs, mirroring
#    how Word would have recorded three distinct edit/typing sessions.
$part1 = "(This is a change " + [char]0x2013 + " Ve"
$rRun1 = $d.Range($pos, $pos)
$rRun1.InsertAfter($part1)
$rRun1.Font.Color = 255
$pos = $pos + $part1.Length

$part2 = "rsion for main branch"
$rRun2 = $d.Range($pos, $pos)
$rRun2.InsertAfter($part2)
$rRun2.Font.Color = 255
$pos = $pos + $part2.Length

$part3 = ")"
$rRun3 = $d.Range($pos, $pos)
$rRun3.InsertAfter($part3)
$rRun3.Font.Color = 255
